$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Last update" timestamp string (shared string, single reference at A300)
$ws.Range("A300").Value = "Last update: 25-02-2025, 10:26"

# Update active selection to A4 (matches sheetView selection in diff)
$ws.Range("A4").Select() | Out-Null

# --- Data cell updates ---
$ws.Range("G13").Copy()
$ws.Range("T13").PasteSpecial(-4122)
$ws.Range("T13").Value = 0.05

$ws.Range("R15").ClearContents()
$ws.Range("E15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("G15").Copy()
$ws.Range("T15").PasteSpecial(-4122)
$ws.Range("T15").Value = 3.1

$ws.Range("P46").Copy()
$ws.Range("T46").PasteSpecial(-4122)
$ws.Range("T46").Value = 37.6

$ws.Range("P47").Copy()
$ws.Range("T47").PasteSpecial(-4122)
$ws.Range("T47").Value = 65

$ws.Range("G61").Copy()
$ws.Range("T61").PasteSpecial(-4122)
$ws.Range("T61").Value = 99.9

$ws.Range("G61").Copy()
$ws.Range("U61").PasteSpecial(-4122)
$ws.Range("U61").Value = 99.9

$ws.Range("R96").Copy()
$ws.Range("U96").PasteSpecial(-4122)
$ws.Range("U96").Value = 22

$ws.Range("R97").Copy()
$ws.Range("U97").PasteSpecial(-4122)
$ws.Range("U97").Value = 20

$ws.Range("R98").Copy()
$ws.Range("U98").PasteSpecial(-4122)
$ws.Range("U98").Value = 25

$ws.Range("R99").Value = 25

$ws.Range("R99").Copy()
$ws.Range("U99").PasteSpecial(-4122)
$ws.Range("U99").Value = 27

$ws.Range("R100").Value = 14

$ws.Range("R100").Copy()
$ws.Range("U100").PasteSpecial(-4122)
$ws.Range("U100").Value = 15

$ws.Range("R101").Copy()
$ws.Range("U101").PasteSpecial(-4122)
$ws.Range("U101").Value = 36

$ws.Range("R102").Copy()
$ws.Range("U102").PasteSpecial(-4122)
$ws.Range("U102").Value = 37

$ws.Range("R103").Copy()
$ws.Range("U103").PasteSpecial(-4122)
$ws.Range("U103").Value = 31

$ws.Range("R104").Copy()
$ws.Range("U104").PasteSpecial(-4122)
$ws.Range("U104").Value = 19

$ws.Range("R105").Copy()
$ws.Range("U105").PasteSpecial(-4122)
$ws.Range("U105").Value = 12

$ws.Range("R106").Copy()
$ws.Range("U106").PasteSpecial(-4122)
$ws.Range("U106").Value = 3

$ws.Range("R107").Copy()
$ws.Range("U107").PasteSpecial(-4122)
$ws.Range("U107").Value = 18

$ws.Range("R108").Copy()
$ws.Range("U108").PasteSpecial(-4122)
$ws.Range("U108").Value = 13

$ws.Range("R109").Copy()
$ws.Range("U109").PasteSpecial(-4122)
$ws.Range("U109").Value = 44

$ws.Range("R110").Copy()
$ws.Range("U110").PasteSpecial(-4122)
$ws.Range("U110").Value = 29

$ws.Range("R111").Copy()
$ws.Range("U111").PasteSpecial(-4122)
$ws.Range("U111").Value = 30

$ws.Range("R112").Copy()
$ws.Range("U112").PasteSpecial(-4122)
$ws.Range("U112").Value = 12

$ws.Range("R113").Copy()
$ws.Range("U113").PasteSpecial(-4122)
$ws.Range("U113").Value = 11

$ws.Range("G114").Copy()
$ws.Range("T114").PasteSpecial(-4122)
$ws.Range("T114").Value = 7.9

$ws.Range("G132").Copy()
$ws.Range("T132").PasteSpecial(-4122)
$ws.Range("T132").Value = 0.2

$ws.Range("R138").Value = 15.62

$ws.Range("S138").Value = 16.89

$ws.Range("G138").Copy()
$ws.Range("T138").PasteSpecial(-4122)
$ws.Range("T138").Value = 16.559999999999999

$ws.Range("G141").Copy()
$ws.Range("T141").PasteSpecial(-4122)
$ws.Range("T141").Value = 9.5

$ws.Range("Q191").Value = 12.1

$ws.Range("R191").Value = 11.9

$ws.Range("S191").Value = 13.1

$ws.Range("G191").Copy()
$ws.Range("T191").PasteSpecial(-4122)
$ws.Range("T191").Value = 13.3

$ws.Range("G192").Copy()
$ws.Range("T192").PasteSpecial(-4122)
$ws.Range("T192").Value = 83.3

$ws.Range("G194").Copy()
$ws.Range("T194").PasteSpecial(-4122)
$ws.Range("T194").Value = 27

$ws.Range("J195").Copy()
$ws.Range("T195").PasteSpecial(-4122)
$ws.Range("T195").Value = 26.7

$ws.Range("J196").Copy()
$ws.Range("T196").PasteSpecial(-4122)
$ws.Range("T196").Value = 26.6

$ws.Range("G208").Copy()
$ws.Range("T208").PasteSpecial(-4122)
$ws.Range("T208").Value = 21.6

$ws.Range("G209").Copy()
$ws.Range("T209").PasteSpecial(-4122)
$ws.Range("T209").Value = 144.6

$ws.Range("G214").Copy()
$ws.Range("T214").PasteSpecial(-4122)
$ws.Range("T214").Value = 94.8

$ws.Range("G215").Copy()
$ws.Range("T215").PasteSpecial(-4122)
$ws.Range("T215").Value = 96.8

$ws.Range("Q217").Copy()
$ws.Range("T217").PasteSpecial(-4122)
$ws.Range("T217").Value = 13.3

$ws.Range("G223").Copy()
$ws.Range("T223").PasteSpecial(-4122)
$ws.Range("T223").Value = 144.6

$ws.Range("S231").Value = 0.97

$ws.Range("G231").Copy()
$ws.Range("T231").PasteSpecial(-4122)
$ws.Range("T231").Value = 1.1399999999999999

$ws.Range("R232").Value = 18.25

$ws.Range("S232").Value = 18.309999999999999

$ws.Range("G232").Copy()
$ws.Range("T232").PasteSpecial(-4122)
$ws.Range("T232").Value = 18.03

$ws.Range("Q233").Value = 7.4

$ws.Range("R233").Value = 7

$ws.Range("S233").Value = 6.7

$ws.Range("G233").Copy()
$ws.Range("T233").PasteSpecial(-4122)
$ws.Range("T233").Value = 7.5

$ws.Range("R234").ClearContents()
$ws.Range("E234").Copy()
$ws.Range("R234").PasteSpecial(-4122)

$ws.Range("G234").Copy()
$ws.Range("T234").PasteSpecial(-4122)
$ws.Range("T234").Value = 3.1

$ws.Range("R237").Value = 15.62

$ws.Range("S237").Value = 16.89

$ws.Range("G237").Copy()
$ws.Range("T237").PasteSpecial(-4122)
$ws.Range("T237").Value = 16.559999999999999

$ws.Range("G238").Copy()
$ws.Range("T238").PasteSpecial(-4122)
$ws.Range("T238").Value = 1372

$ws.Range("G242").Copy()
$ws.Range("S242").PasteSpecial(-4122)
$ws.Range("S242").Value = 1959

$ws.Range("G243").Value = 48641

$ws.Range("L243").Value = 46834

$ws.Range("Q243").Value = 61823

$ws.Range("G243").Copy()
$ws.Range("R243").PasteSpecial(-4122)
$ws.Range("R243").Value = 66958

$ws.Range("G243").Copy()
$ws.Range("S243").PasteSpecial(-4122)
$ws.Range("S243").Value = 68339

$ws.Range("G254").Copy()
$ws.Range("T254").PasteSpecial(-4122)
$ws.Range("T254").Value = 205.6

$ws.Range("G257").Copy()
$ws.Range("T257").PasteSpecial(-4122)
$ws.Range("T257").Value = 0.192

$ws.Range("G264").Copy()
$ws.Range("T264").PasteSpecial(-4122)
$ws.Range("T264").Value = 123796

$ws.Range("N268").Value = 0.81

$ws.Range("G268").Copy()
$ws.Range("T268").PasteSpecial(-4122)
$ws.Range("T268").Value = 0.78

$ws.Range("T289").Value = 2580.38

$ws.Range("T290").Value = 1701.94

$ws.Range("T291").Value = 0.33

$ws.Range("T293").Value = 4.2699999999999996

$ws.Range("T297").Copy()
$ws.Range("U297").PasteSpecial(-4122)
$ws.Range("U297").Value = 61
